$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "https://maps.app.goo.gl/1WURXYDUiB8qZMct9"
$ws.Range("C7").Value = "Masjid Abu Bakr As-Siddiq"
$ws.Range("E7").Value = "https://maps.app.goo.gl/UcVxLZGibCLYtGPa8"
$ws.Range("C14").Value = "Masjid Umar bin Khattab"
$ws.Range("E14").Value = "https://maps.app.goo.gl/51CtZM5THJfdUoYH8"
$ws.Range("E21").Value = "https://maps.app.goo.gl/d6Bom9oBgSpAZ6Ag8"
$ws.Range("C37").Value = "Bir-e-Shifa Well"
$ws.Range("E37").Value = "https://maps.app.goo.gl/cdBW3f3ozkDPgQHz5"
$ws.Range("E55").Value = "https://maps.app.goo.gl/k4DhZg1Wdsq3Zmik8"
$ws.Range("C70").Value = "Masjid al-Fateh"
$ws.Range("C72").Value = "Masjid Ali bin Abi Talib"
$ws.Range("E72").Value = "https://maps.app.goo.gl/s5j5AtRQhoPK1SGP8"
$ws.Range("B73").Value = "places"
$ws.Range("C73").Value = "Masjid Salman-al-Farsi"
$ws.Range("D73").Value = "Salman Al Farsi Mosque"
$ws.Range("E73").Value = "https://maps.app.goo.gl/gJu3KAnetLyR7q229"
$ws.Range("C74").Value = "Fastest Arrow Dates"
$ws.Range("D74").Value = "Premium dates, nuts, honey and chocolate. Very popular for premium real honey and high quality Ajwa dates"
$ws.Range("A75").Value = "MADINAH"
$ws.Range("C75").Value = "Turkish Station Underpass"
$ws.Range("D75").Value = "Cheap local goods."
$ws.Range("C76").Value = "Al Marwah Shop"
$ws.Range("D76").Value = "Cheapest gift shop in Makkah."
$ws.Range("C77").Value = "Bin Dawood (Aziziyah Branch)"
$ws.Range("D77").Value = "Supermarket chain for groceries, clothes & souvenirs. Groceries & souvenirs, 3-4 km from Haram."
$ws.Range("A78").Value = "MAKKAH"
$ws.Range("B78").Value = "shopping"
$ws.Range("C78").Value = "Top 10 Shop (Aziziyah)"
$ws.Range("D78").Value = "Everything around 10 SAR (affordable gifts, toys, home items). Budget items (~10 SAR)."
$ws.Range("B79").Value = "food"
$ws.Range("C79").Value = "Gewar Taiba"
$ws.Range("D79").Value = "Shawarma."
$ws.Range("A80").Value = "MADINAH"
$ws.Range("B80").Value = "places"
$ws.Range("C80").Value = "The International Fair & Museum of the Prophet's Biography and Islamic Civilization"
$ws.Range("D80").Value = "Near Masjid an-Nabawi, immersive museum with VR, holography, and interactive displays."
$ws.Range("E80").Value = "https://maps.app.goo.gl/YdeRo66f5MthvDSU9"
$ws.Range("C81").Value = "Mama Ghazzel"
$ws.Range("D81").Value = "Sweets & desserts, Jabal E Omar area."
$ws.Range("C82").Value = "Lamma Burger Crispy Roll / Red Burger"
$ws.Range("D82").Value = "Food to Try."
$ws.Range("B83").Value = "food"
$ws.Range("C83").Value = "Nimra Shinwari Hotel"
$ws.Range("D83").Value = "Peshawari Cuisine (Aziziya). Peshawari food, Aziziya."
$ws.Range("B84").Value = "shopping"
$ws.Range("C84").Value = "Abayas"
$ws.Range("D84").Value = "Ajyad, Misfalah & Al Diyafa streets - SAR 40 to 100."
$ws.Range("B85").Value = "food"
$ws.Range("C85").Value = "Ice Cream Al Asemah"
$ws.Range("D85").Value = "Off Al Ghufran Hotel, Ajyad Street (try Orange Slush, 10 SAR)."
$ws.Range("A86").Value = "MAKKAH"
$ws.Range("B86").Value = "places"
$ws.Range("C86").Value = "Masjid al-Jinn"
$ws.Range("D86").Value = "Distance: 1-3 km. Specialized for: Prophet recited Qur'an to jinn. Tip: Short visit, accessible by foot."
$ws.Range("B87").Value = "shopping"
$ws.Range("C87").Value = "Dawoodiya Market (Quba Road, Al Jumuah)"
$ws.Range("D87").Value = "Abayas from SAR 30 upwards."
$ws.Range("C88").Value = "Uhud Mountain & Martyrs' Cemetery"
$ws.Range("D88").Value = "Shuhada Uhud Cemetery`n6-7 km."
$ws.Range("E88").Value = "https://maps.app.goo.gl/T7RYXXY66svqEQt1A"
$ws.Range("A89").Value = "MADINAH"
$ws.Range("C89").Value = "Masjid al-Bilal (Bilal Ibn Rabah Mosque)"
$ws.Range("D89").Value = "Near Quba, named after first muezzin. 4-6 km."
$ws.Range("E89").Value = "https://maps.app.goo.gl/67eRWTuV7jhwLa5J9"
$ws.Range("B90").Value = "places"
$ws.Range("C90").Value = "Taif Day Trip"
$ws.Range("D90").Value = "Distance: 90-95 km. Specialized for: Rose gardens & markets, Rose distilleries & mountain views. Tip: Day trip, best in morning."
$ws.Range("C91").Value = "Barn's Coffee"
$ws.Range("D91").Value = "Local Saudi coffee chain. Coffee & snacks."
$ws.Range("B92").Value = "food"
$ws.Range("C92").Value = "Almarai 100% Apple / Mixed Fruit Juice"
$ws.Range("D92").Value = "(no sugar) Food to Try."
$ws.Range("B93").Value = "shopping"
$ws.Range("C93").Value = "Abraj Hypermarket"
$ws.Range("D93").Value = "Groceries & souvenirs. General shopping."
$ws.Range("A94").Value = "MAKKAH"
$ws.Range("B94").Value = "food"
$ws.Range("C94").Value = "Funduq Undlusia Restaurant"
$ws.Range("D94").Value = "Known for Arabic set meals. Local cuisine restaurant."
$ws.Range("B95").Value = "shopping"
$ws.Range("C95").Value = "Premium Janimaz Shop"
$ws.Range("D95").Value = "Gate 305, prayer mats."
$ws.Range("B96").Value = "places"
$ws.Range("C96").Value = "Masjid Quba"
$ws.Range("D96").Value = "4-6 km, First mosque in Islam."
$ws.Range("E96").Value = "https://maps.app.goo.gl/63AKEaKmLk7LKbyw5"
$ws.Range("A97").Value = "MADINAH"
$ws.Range("B97").Value = "food"
$ws.Range("C97").Value = "Karak Express"
$ws.Range("D97").Value = "Chicken Biryani & Karak Chai. Chai & Chicken Biryani."
$ws.Range("B98").Value = "shopping"
$ws.Range("C98").Value = "Shamali Aziziya Mobile Market"
$ws.Range("D98").Value = "Cheap iPhones."
$ws.Range("A99").Value = "MAKKAH"
$ws.Range("B99").Value = "food"
$ws.Range("C99").Value = "Dunkin Donuts Cold Brew"
$ws.Range("D99").Value = "Food to Try."
$ws.Range("A100").Value = "MADINAH"
$ws.Range("C100").Value = "Awaali Gardens"
$ws.Range("D100").Value = "Historic palm groves gifted to Prophet."
$ws.Range("A101").Value = "MAKKAH"
$ws.Range("B101").Value = "places"
$ws.Range("C101").Value = "Masjid al-Ji'rana"
$ws.Range("D101").Value = "Distance: 24-26 km. Specialized for: Miqat for locals / Miqat for residents of Makkah for Umrah, historic site after Battle of Hunayn."
$ws.Range("A102").Value = "MADINAH"
$ws.Range("B102").Value = "shopping"
$ws.Range("C102").Value = "Tamil Nadu Dates Shop"
$ws.Range("D102").Value = "Shop No. 9, Near Gate 330 (Taiba Center side)."
$ws.Range("A103").Value = "MAKKAH"
$ws.Range("C103").Value = "Nadec 100% Apple / Strawberry Juice"
$ws.Range("D103").Value = "Food to Try."
$ws.Range("B104").Value = "food"
$ws.Range("C104").Value = "Qasr-ul-Amal Restaurant"
$ws.Range("D104").Value = "Traditional Arabic meals."
$ws.Range("A105").Value = "MADINAH"
$ws.Range("B105").Value = "places"
$ws.Range("C105").Value = "Masjid al-Khandaq (The Seven Mosques)"
$ws.Range("E105").Value = "https://maps.app.goo.gl/v5VCu6RXfo1xmvQk9"
